# Generate Report for Handback
# - Marks the b870009e... row's Status as "Handback transform failed"
#   (instead of "Ready for handoff") on the Overview, zh-cn and de-de sheets.
# - Records the handback-filename-mismatch error detail on row 3
#   (Error Detail column, L) of both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status columns for the b870009e row (row 3) ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Handback transform failed"
$ovw.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("L3").Value = "Handback file name: c1yse3t2.oo0 is different with handoff file name: b870009e-7fc9-4f00-a3b4-e455a10b01bb.12e88d09951fe179cf3f3d3c4ee040e8ef09034a.zh-cn."

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("L3").Value = "Handback file name: c1yse3t2.oo0 is different with handoff file name: b870009e-7fc9-4f00-a3b4-e455a10b01bb.12e88d09951fe179cf3f3d3c4ee040e8ef09034a.de-de."
